# Generate Report for Handoff
# Updates the "Latest Handoff" timestamp columns for the row corresponding to
# 5759e26c-d66c-4147-8177-feafb869c911.md on all three report sheets, as a
# new handoff xliff generation was recorded.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
# 5759e26c-d66c-4147-8177-feafb869c911.md row (row 4).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-09-05 10:08:58"

# zh-cn sheet: "Latest Handoff Datetime" column (H) for the same file (row 4).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-09-05 10:08:48"

# de-de sheet: "Latest Handoff Datetime" column (H) for the same file (row 4).
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-09-05 10:08:58"
